# Weekly price-sheet update: insert a new daily record for "Mandarina"
# (Murcott / Primera) right before the existing row 170, pushing every
# subsequent record down by one row (170-202 -> 171-203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 170; Excel shifts rows 170..202 down to 171..203
# and extends the used range to A1:T203.
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(170, 1).Value = 4
$ws.Cells.Item(170, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(170, 3).Value = "Los Lagos"
$ws.Cells.Item(170, 4).Value = 44694
$ws.Cells.Item(170, 5).Value = 10
$ws.Cells.Item(170, 6).Value = "Fruta"
$ws.Cells.Item(170, 7).Value = 100102
$ws.Cells.Item(170, 8).Value = "Cítricos"
$ws.Cells.Item(170, 9).Value = 100102004
$ws.Cells.Item(170, 10).Value = "Mandarina"
$ws.Cells.Item(170, 11).Value = "Murcott"
$ws.Cells.Item(170, 12).Value = "Primera"
$ws.Cells.Item(170, 13).Value = 700
$ws.Cells.Item(170, 14).Value = 14000
$ws.Cells.Item(170, 15).Value = 15000
$ws.Cells.Item(170, 16).Value = 14500
$ws.Cells.Item(170, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(170, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(170, 19).Value = 1450
$ws.Cells.Item(170, 20).Value = 10
